$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# New handback report run:
#   - the previously-handed-back file "5618eef6-2572-4309-abff-b8645fe5ce31"
#     was re-run and is now reported under a new run id
#     "6c277fea-61dd-4d0d-bf6d-5662fc177054" with refreshed timestamps
#   - a brand new file "7bfec0f9-571c-4e5a-8cf7-dd6f80875346" was handed
#     back in the same run, and needs a new row on every sheet
# ---------------------------------------------------------------------------

$oldId = "5618eef6-2572-4309-abff-b8645fe5ce31"
$newId = "6c277fea-61dd-4d0d-bf6d-5662fc177054"
$addId = "7bfec0f9-571c-4e5a-8cf7-dd6f80875346"

$oldZhXlf = "$oldId.eedaa5a4e5c001da6d6e901a393d35f43e221077.zh-cn.xlf"
$oldDeXlf = "$oldId.eedaa5a4e5c001da6d6e901a393d35f43e221077.de-de.xlf"

$newZhXlf = "$newId.bd094a099b0b2be849868b7ff0281bc98b158a94.zh-cn.xlf"
$newDeXlf = "$newId.bd094a099b0b2be849868b7ff0281bc98b158a94.de-de.xlf"

$addZhXlf = "$addId.9f71b44250bb45e9c725ae1272cb1547bcbfc922.zh-cn.xlf"
$addDeXlf = "$addId.9f71b44250bb45e9c725ae1272cb1547bcbfc922.de-de.xlf"

$latestHoDate   = "2016-08-15 16:57:49"
$zhHandoffDate  = "2016-08-15 16:57:44"
$zhHandbackDate = "2016-08-15 16:58:05"
$deHandbackDate = "2016-08-15 16:58:15"

# ---------------------------------------------------------------------------
# Grow each table by one row so the new file has a place to live.
# ---------------------------------------------------------------------------
$wsOverview.ListObjects.Item(1).ListRows.Add() | Out-Null
$wsZhCn.ListObjects.Item(1).ListRows.Add() | Out-Null
$wsDeDe.ListObjects.Item(1).ListRows.Add() | Out-Null

# ---------------------------------------------------------------------------
# Overview sheet: row 2 is the renamed run, row 3 is the new file.
# Columns: A File Name, B Path And Name, C Extension, D Publish URL,
#          E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = $latestHoDate

$wsOverview.Range("A3").Value = "$addId.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = $latestHoDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/501147d9e309f4b6c7e21dee3849dd8e7fd4d96f/e2e/$addId.md", `
    "", "", "e2e\$addId.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: row 2 is the renamed run, row 3 is the new file.
# Columns: A Source File Name, B File Extension, C Status, D Source Path,
#          E Priority, F Content Duplicate, G Correspond Handoff File,
#          H Correspond Handoff Datetime, I Target File,
#          J Correspond Handback File, K Correspond Handback DateTime,
#          L Reference Tokens, M To be localized, N Dependency From,
#          O Has metadata, P Error Detail
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $zhHandoffDate
$wsZhCn.Range("I2").Value = "$newId.md"
$wsZhCn.Range("J2").Value = $newZhXlf
$wsZhCn.Range("K2").Value = $zhHandbackDate

$wsZhCn.Range("A3").Value = "$addId.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = $addZhXlf
$wsZhCn.Range("H3").Value = $zhHandoffDate
$wsZhCn.Range("I3").Value = "$addId.md"
$wsZhCn.Range("J3").Value = $addZhXlf
$wsZhCn.Range("K3").Value = $zhHandbackDate
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/501147d9e309f4b6c7e21dee3849dd8e7fd4d96f/e2e/$addId.md", `
    "", "", "$addId.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5d698a1a7591d36443a41d22b0093b86dfd18fac/e2e/$addId.md", `
    "", "", "$addId.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: same column layout as zh-cn.
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $latestHoDate
$wsDeDe.Range("I2").Value = "$newId.md"
$wsDeDe.Range("J2").Value = $newDeXlf
$wsDeDe.Range("K2").Value = $deHandbackDate

$wsDeDe.Range("A3").Value = "$addId.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = $addDeXlf
$wsDeDe.Range("H3").Value = $latestHoDate
$wsDeDe.Range("I3").Value = "$addId.md"
$wsDeDe.Range("J3").Value = $addDeXlf
$wsDeDe.Range("K3").Value = $deHandbackDate
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/501147d9e309f4b6c7e21dee3849dd8e7fd4d96f/e2e/$addId.md", `
    "", "", "$addId.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8ab6e4add655c890c0bf329992c52d79a657f97d/e2e/$addId.md", `
    "", "", "$addId.md") | Out-Null
